# Adds a new "2022-Q1" quarterly sheet (before the "总计" summary sheet)
# and records its totals as a new top row in the "总计" sheet.

function Set-TextValue($range, $value) {
    # Excel's COM layer auto-detects numeric-looking strings and stores them
    # as numbers. Several columns in this workbook intentionally keep such
    # values as text (e.g. fund codes, "0.3800" needs its trailing zero).
    # Temporarily formatting the cell as Text forces string storage; we then
    # restore the cell style to the workbook's default so no stray
    # number-format carries over.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert the new "2022-Q1" sheet directly before "总计".
# ------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# Clone the layout/styles of an existing quarterly sheet (header row style +
# the "index" column style) so the new sheet matches the others, then
# overwrite with the real 2022-Q1 figures.
$srcSheet = $wb.Worksheets.Item("2021-Q3")
$srcSheet.Range("A1:H4").Copy($newSheet.Range("A1"))
$newSheet.Cells.Item(1, 1).ClearContents()

# Header row (column D label differs from the older "基金金额" wording).
$newSheet.Cells.Item(1, 2).Value = "基金代码"
$newSheet.Cells.Item(1, 3).Value = "基金名称"
$newSheet.Cells.Item(1, 4).Value = "基金规模"
$newSheet.Cells.Item(1, 5).Value = "股票总仓位"
$newSheet.Cells.Item(1, 6).Value = "仓位占比"
$newSheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1, 8).Value = "仓位排名"

# Row 2 - 广发中证全指电力ETF
$newSheet.Cells.Item(2, 1).Value = 0
Set-TextValue $newSheet.Cells.Item(2, 2) "159611"
$newSheet.Cells.Item(2, 3).Value = "广发中证全指电力ETF"
Set-TextValue $newSheet.Cells.Item(2, 4) "13.38"
Set-TextValue $newSheet.Cells.Item(2, 5) "99.14"
Set-TextValue $newSheet.Cells.Item(2, 6) "2.84"
Set-TextValue $newSheet.Cells.Item(2, 7) "0.3800"
$newSheet.Cells.Item(2, 8).Value = 8

# Row 3 - 银华巨潮小盘价值ETF
$newSheet.Cells.Item(3, 1).Value = 1
Set-TextValue $newSheet.Cells.Item(3, 2) "159990"
$newSheet.Cells.Item(3, 3).Value = "银华巨潮小盘价值ETF"
Set-TextValue $newSheet.Cells.Item(3, 4) "1.06"
Set-TextValue $newSheet.Cells.Item(3, 5) "96.39"
Set-TextValue $newSheet.Cells.Item(3, 6) "1.25"
Set-TextValue $newSheet.Cells.Item(3, 7) "0.0132"
$newSheet.Cells.Item(3, 8).Value = 7

# Row 4 - 中银顺盈回报一年持有期混合
$newSheet.Cells.Item(4, 1).Value = 2
Set-TextValue $newSheet.Cells.Item(4, 2) "010487"
$newSheet.Cells.Item(4, 3).Value = "中银顺盈回报一年持有期混合"
Set-TextValue $newSheet.Cells.Item(4, 4) "1.23"
Set-TextValue $newSheet.Cells.Item(4, 5) "38.83"
Set-TextValue $newSheet.Cells.Item(4, 6) "0.95"
Set-TextValue $newSheet.Cells.Item(4, 7) "0.0117"
$newSheet.Cells.Item(4, 8).Value = 7

# ------------------------------------------------------------------
# 2. Add the 2022-Q1 summary row to the top of "总计"'s data rows.
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("总计")
$ws.Rows(2).Insert()

# The inserted row copies formatting from the row above (the bold header);
# reset it, then restore just the index-column style from a data row below.
$ws.Range("A2:D2").ClearFormats()
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122) # xlPasteFormats

$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "2022-Q1"
$ws.Cells.Item(2, 3).Value = 3
$ws.Cells.Item(2, 4).Value = 0.4

# The shifted-down rows keep their old index-column numbers; renumber them
# (0,1,2,3,4) now that the new row leads the list.
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(6, 1).Value = 4
